$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy per-row formatting from column L into the new column M ---
# (the style used in each M cell mirrors the style already used by the
#  corresponding L cell on the same row, reusing the existing cellXfs entries)

$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)

$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("M13").PasteSpecial(-4122)
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("M21").PasteSpecial(-4122)
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("M24").PasteSpecial(-4122)
$ws.Range("M25").PasteSpecial(-4122)
$ws.Range("M27").PasteSpecial(-4122)
$ws.Range("M28").PasteSpecial(-4122)
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("M31").PasteSpecial(-4122)

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("M32").PasteSpecial(-4122)

$ws.Range("L23").Copy()
$ws.Range("M23").PasteSpecial(-4122)
$ws.Range("M29").PasteSpecial(-4122)

$ws.Range("L33").Copy()
$ws.Range("M33").PasteSpecial(-4122)

# M26 needs a style that does not exist yet: same number format/font as style
# "13" (0.0, 9pt Times New Roman) but right-aligned -- Excel will add a new cellXf
# for this combination, just like in the authored workbook.
$ws.Range("L5").Copy()
$ws.Range("M26").PasteSpecial(-4122)
$ws.Range("M26").HorizontalAlignment = -4152

$excel.CutCopyMode = $false

# --- Set the values for the new column M (2021 data) ---
$ws.Range("M3").Value() = 2021
$ws.Range("M4").Value() = 2.0173148373954581
$ws.Range("M5").Value() = 0.11867182493532386
$ws.Range("M6").Value() = 3.9440914499323179
$ws.Range("M7").Value() = 0
$ws.Range("M8").Value() = "-"
$ws.Range("M9").Value() = 0
$ws.Range("M10").Value() = 0.62921030174566528
$ws.Range("M11").Value() = "-"
$ws.Range("M12").Value() = 1.2497227177719943
$ws.Range("M13").Value() = 0.19844537890168421
$ws.Range("M14").Value() = "-"
$ws.Range("M15").Value() = 0.39861918314956984
$ws.Range("M16").Value() = 0
$ws.Range("M17").Value() = "-"
$ws.Range("M18").Value() = 0
$ws.Range("M19").Value() = 0.85521252031129735
$ws.Range("M20").Value() = "-"
$ws.Range("M21").Value() = 1.6913581464969858
$ws.Range("M22").Value() = 1.8347815875998121
$ws.Range("M23").Value() = "-"
$ws.Range("M24").Value() = 3.6321107648498847
$ws.Range("M25").Value() = 6.1211560415300026
$ws.Range("M26").Value() = "-"
$ws.Range("M27").Value() = 12.437939862560766
$ws.Range("M28").Value() = 3.6823562661275693
$ws.Range("M29").Value() = 0.69433233870225819
$ws.Range("M30").Value() = 7.0564990356117976
$ws.Range("M31").Value() = 2.7447727328177227
$ws.Range("M32").Value() = "-"
$ws.Range("M33").Value() = 5.6418550419377889

# --- Restore the active selection, matching the authored workbook ---
$ws.Range("P6").Select()
